# Daily attendance processing - 2026-01-01 20:35:45
# Reverses the order of the comma-separated "Recorded By" entries in column G
# for every data row on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value()

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[-1..-$parts.Count]
            $newVal = [string]::Join(", ", $reversed)
            $cell.Value = $newVal
        }
    }
}
